$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet (3rd sheet) - insert a new blank
# column before column N, shifting the "Late"/"heading"/"Outstanding"
# columns (and their data) one place to the right.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N:N").Insert()

# The inserted column picks up the width of its left neighbour (column M),
# matching Excel's native "insert column" behaviour.
$wsRepay.Columns("N:N").ColumnWidth = $wsRepay.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, and update its selection.
# (This also clears the previous tabSelected flag that was on "Input".)
$wsRepay.Activate()
$wsRepay.Range("J16").Select()
